$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reassign the "ward map" related rows: move several tasks that were
# further down the sheet up into a new contiguous block starting at row 25,
# and update the entry that used to live at row 37 (now row 47) so its
# urgency changes from "high" to "low".

# First, clear out the old source rows entirely so they disappear from the
# sheet once their data has been relocated.
$ws.Range("A36:E36").ClearContents()
$ws.Range("A37:E37").ClearContents()
$ws.Range("A39:E39").ClearContents()
$ws.Range("A40:E40").ClearContents()
$ws.Range("A46:E46").ClearContents()
$ws.Range("A47:E47").ClearContents()
$ws.Range("A48:E48").ClearContents()
$ws.Range("A49:E49").ClearContents()
$ws.Range("A57:E57").ClearContents()

# The destination rows (25, 27-32) were blank placeholder rows that only
# carried a style on B/C -- clear that leftover formatting since the new
# rows (other than 29) should have no explicit style.
$ws.Range("B25:C25").ClearFormats()
$ws.Range("B27:C28").ClearFormats()
$ws.Range("B30:C32").ClearFormats()

# Row 25 <- old row 36 ("add distro chart x axis label"), now also marked done
$ws.Range("A25").Value = "add distro chart x axis label"
$ws.Range("B25").Value = "charts"
$ws.Range("C25").Value = "high"
$ws.Range("D25").Value = "done"

# Row 27 <- old row 46 ("make dots darker in maps?")
$ws.Range("A27").Value = "make dots darker in maps?"
$ws.Range("B27").Value = "maps"
$ws.Range("C27").Value = "medium"
$ws.Range("E27").Value = "Ask Mal????"

# Row 28 <- old row 49 ("apply copy edits")
$ws.Range("A28").Value = "apply copy edits"
$ws.Range("B28").Value = "body"
$ws.Range("C28").Value = "high"

# Row 29 <- old row 40 ("make charts load better on mobile") - keeps its style
$ws.Range("A29").Value = "make charts load better on mobile"
$ws.Range("B29").Value = "charts"
$ws.Range("C29").Value = "high"

# Row 30 <- old row 47 ("Improve dot fade experience")
$ws.Range("A30").Value = "Improve dot fade experience"
$ws.Range("B30").Value = "body"
$ws.Range("C30").Value = "high"

# Row 31 <- old row 48 ("dot fade on mobile etc")
$ws.Range("A31").Value = "dot fade on mobile etc"
$ws.Range("B31").Value = "body"
$ws.Range("C31").Value = "high"

# Row 32 <- old row 39 ("full responsive on load")
$ws.Range("A32").Value = "full responsive on load"
$ws.Range("B32").Value = "all"
$ws.Range("C32").Value = "high"

# Row 47 (new content) <- old row 37 ("fix scale on ward distros"), urgency
# changed from "high" to "low"
$ws.Range("A47").Value = "fix scale on ward distros"
$ws.Range("B47").Value = "charts"
$ws.Range("C47").Value = "low"
$ws.Range("E47").Value = "make it not overlap? Ask ben"

# Row 56 <- old row 57 ("update social media buttons") shifts up one row
$ws.Range("A56").Value = "update social media buttons"
$ws.Range("B56").Value = "body"
$ws.Range("C56").Value = "high"
$ws.Range("E56").Value = "talk to team"

# --- Update the saved view state to match where the user ended up editing.
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("A27:E27").Select()
